# EPBDS-7754 Rename date() to Date()
# Date() is more constructor-like and more friendly for understanding for BAs
#
# Two edits on the "date" worksheet of the Dates.xlsx test workbook:
#   1. Rename the worksheet tab from "date" to "Date".
#   2. Rename the sample rule body text from "return date(year, month, day);"
#      to "return Date(year, month, day);" (cell B5, merged B5:D5).

$wb = $excel.ActiveWorkbook

# 1. Rename the "date" sheet to "Date".
$ws = $wb.Worksheets.Item("date")
$ws.Name = "Date"

# 2. Update the sample code snippet to call the renamed Date() constructor.
$ws.Range("B5").Value = "return Date(year, month, day);"
